# Add data for 2023-06-01
# Updates YTD crime-count figures on the citywide/neighborhood summary
# sheets plus the individual neighborhood sheets affected by the new day
# of data (two new cells are introduced where a neighborhood had no prior
# value for that crime category/year).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 56
$ws.Range("B3").Value = 33
$ws.Range("I3").Value = 79
$ws.Range("B6").Value = 160
$ws.Range("C6").Value = 203
$ws.Range("D6").Value = 178
$ws.Range("F6").Value = 210
$ws.Range("H6").Value = 162
$ws.Range("J6").Value = 173
$ws.Range("B7").Value = 215
$ws.Range("C7").Value = 271
$ws.Range("D7").Value = 275
$ws.Range("F7").Value = 297
$ws.Range("H7").Value = 252
$ws.Range("I7").Value = 354
$ws.Range("J7").Value = 315

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B6").Value = 10
$ws.Range("B13").Value = 2
$ws.Range("B17").Value = 3
$ws.Range("I19").Value = 8
$ws.Range("B26").Value = 17
$ws.Range("C26").Value = 19
$ws.Range("J30").Value = 11
$ws.Range("B51").Value = 17
$ws.Range("I51").Value = 56
$ws.Range("J59").Value = 2
$ws.Range("I60").Value = 7
$ws.Range("H81").Value = 2
$ws.Range("D83").Value = 2
$ws.Range("F92").Value = 4
$ws.Range("B94").Value = 215
$ws.Range("C94").Value = 271
$ws.Range("D94").Value = 275
$ws.Range("F94").Value = 297
$ws.Range("H94").Value = 252
$ws.Range("I94").Value = 354
$ws.Range("J94").Value = 315

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("B5").Value = 9
$ws.Range("B6").Value = 10

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I3").Value = 3
$ws.Range("I5").Value = 8

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 18
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 19

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J6").Value = 11
$ws.Range("J2").Value = 1

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 4

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 3
$ws.Range("B3").Value = 6
$ws.Range("B7").Value = 17
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 2
$ws.Range("I5").Value = 7

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("H5").Value = 2
$ws.Range("H4").Value = 1
